# Update the Sprint2 Plan sheet so that the "Dependency" column (C) records
# which task each item depends on (was previously left blank).
$wb = $excel.ActiveWorkbook
$plan = $wb.Worksheets.Item("Plan")

$plan.Range("C4").Value  = "T:1"
$plan.Range("C5").Value  = "T:2"
$plan.Range("C6").Value  = "T:1"
$plan.Range("C7").Value  = "T:2"
$plan.Range("C8").Value  = "T:1"
$plan.Range("C9").Value  = "T:2"
$plan.Range("C10").Value = "T:2"
$plan.Range("C11").Value = "T:1"
$plan.Range("C12").Value = "T:2"
$plan.Range("C13").Value = "T:1"

# C12/C13 previously carried a redundant "apply fill" flag in their cell
# format; normalize it (no visual effect - fill stays "none") to match the
# formatting used by the rest of the Dependency column.
$plan.Range("C12:C13").Interior.Pattern = -4142

# Restore the active cell/selection that was left on the Plan sheet.
$plan.Activate()
$plan.Range("F18").Select()
